$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the StatQuery text shared by C2:C4 ---
$newQuery = " MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n      WHERE (size([]) = 0 OR s.clinical_study_designation IN [])`n        AND (s.study_disposition = 'Unrestricted')`n        AND (size([]) = 0 OR s.clinical_study_type IN [])`n        AND (size([]) = 0 OR demo.breed IN [])`n        AND (size([]) = 0 OR demo.sex IN [])`n        AND (size([]) = 0 OR demo.neutered_indicator IN [])`n        AND (size([]) = 0 OR diag.disease_term IN [])`n        AND (size([]) = 0 OR diag.primary_disease_site IN [])`n        AND (size([]) = 0 OR diag.stage_of_disease IN [])`n        AND (size([]) = 0 OR diag.best_response IN [])`n    OPTIONAL MATCH (c)-->(co:cohort)`n    OPTIONAL MATCH (f:file)-[*]->(c)`n    OPTIONAL MATCH (f)-->(parent)`n    OPTIONAL MATCH (samp:sample)-->(c)`n    OPTIONAL MATCH (samp)<--(al:aliquot)`n    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al`n      WHERE (size([]) = 0 OR samp.summarized_sample_type IN [])`n        AND (size(['Pulmonary Adenocarcinoma']) = 0 OR samp.specific_sample_pathology IN ['Pulmonary Adenocarcinoma'])`n        AND (size([]) = 0 OR samp.sample_site IN [])`n        AND (size([]) = 0 OR head(labels(parent)) IN [])`n        AND (size([]) = 0 OR f.file_type IN [])`n        AND (size([]) = 0 OR f.file_format IN [])`n    WITH c.case_id AS case_id,`n         s.clinical_study_designation AS study_code,`n         s.clinical_study_type AS study_type,`n         co.cohort_description AS cohort,`n         demo.breed AS breed,`n         diag.disease_term AS diagnosis,`n         diag.stage_of_disease AS stage_of_disease,`n         diag.primary_disease_site AS disease_site,`n         demo.patient_age_at_enrollment AS age,`n         demo.sex AS sex,`n         demo.neutered_indicator AS neutered_status,`n         demo.weight AS weight,`n         diag.best_response AS response_to_treatment,`n         samp.sample_id AS sample_id,`n         f.uuid AS file_id,`n         al`n    RETURN`nCOUNT(DISTINCT file_id) as number_of_files,`nCOUNT(DISTINCT sample_id) as number_of_sample,`nCOUNT(DISTINCT case_id) as number_of_cases,`nCOUNT(DISTINCT study_code) as number_of_study,`nCOUNT(DISTINCT al) as number_of_aliquot`n    "

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# --- Row heights grow to fit the much longer query text (capped at Excel's max) ---
$ws.Rows.Item(2).RowHeight = 409.6
$ws.Rows.Item(3).RowHeight = 409.6
$ws.Rows.Item(4).RowHeight = 409.6

# --- View / selection changes ---
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("C4").Select()
try { $win.ScrollRow = 4 } catch {}
try { $win.ScrollColumn = 1 } catch {}

# --- Best-effort: window geometry (engine may not persist this to XML) ---
try {
    $win.Width = 1162.8
    $win.Height = 628.8
    $win.Top = -5.4
    $win.Left = -5.4
} catch {}
